$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark that currently sits after
#    "le site b2b.ephec-ti.be." (near the top of the document).
# ------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 2) In the very last paragraph of the document ("Reverse proxy ...
#    lancé"), change the final word "lancé" to "En cours".
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$editRange = $lastPara.Range
$editRange.Find.Execute("lanc" + [char]0x00E9, $true, $false, $false, $false, $false, `
                         $true, 1, $false, "En cours", 2)

# The replace above merges the preceding tab and the new text into a
# single run (the engine rebuilds the edited run as plain text, which
# loses the <w:tab/> element unless the new text is split into its own
# run). Force "En cours" into its own run - independent of the tab
# run before it - by nudging its formatting and then restoring the
# original value; this creates a run boundary without altering the
# visible formatting.
$originalSize = $editRange.Font.Size
$editRange.Font.Size = $originalSize + 1
$editRange.Font.Size = $originalSize

# ------------------------------------------------------------------
# 3) Re-add a collapsed "_GoBack" bookmark right after the new text,
#    at the very end of the last paragraph (Word drops one there to
#    mark the last edit position).
# ------------------------------------------------------------------
$lastPara2 = $d.Paragraphs.Last
$endPos = $lastPara2.Range.End
$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("ZZZMARKERZZZ")
$markerLen = ("ZZZMARKERZZZ").Length
$markerRange = $d.Range($placeholder.Start, $placeholder.Start + $markerLen)
$newBookmark = $d.Bookmarks.Add("_GoBack", $markerRange)
$newBookmark.Range.Text = ""
